# Update gh-pages to output generated at 456a3b4
# Updates the F-column "视频数/指标" values on the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7478
$ws1.Range("F3").Value = 69
$ws1.Range("F4").Value = 212
$ws1.Range("F5").Value = 229
$ws1.Range("F6").Value = 1125
$ws1.Range("F8").Value = 19
$ws1.Range("F9").Value = 122
$ws1.Range("F10").Value = 31

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7478
$ws4.Range("F3").Value = 69
$ws4.Range("F4").Value = 212
$ws4.Range("F5").Value = 229
$ws4.Range("F6").Value = 1125
$ws4.Range("F9").Value = 19
$ws4.Range("F10").Value = 122
$ws4.Range("F11").Value = 31
